$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '29.391.82'
$ws.Range('E2').Value = '  +9.90%  '
Set-TextValue 'D3' '1.837.80'
$ws.Range('E3').Value = '  +6.87%  '
Set-TextValue 'D4' '0.9982'
$ws.Range('E4').Value = '  -0.23%  '
Set-TextValue 'D5' '247.09'
$ws.Range('E5').Value = '  +3.09%  '
Set-TextValue 'D6' '0.9978'
$ws.Range('E6').Value = '  -0.30%  '
Set-TextValue 'D7' '0.4938'
$ws.Range('E7').Value = '  +3.96%  '
Set-TextValue 'D8' '0.2803'
$ws.Range('E8').Value = '  +9.95%  '
$ws.Range('E9').Value = '  +4.95%  '
Set-TextValue 'D10' '1.825.44'
$ws.Range('E10').Value = '  +6.23%  '
Set-TextValue 'D11' '16.80'
$ws.Range('E11').Value = '  +6.23%  '
Set-TextValue 'D12' '0.07109'
$ws.Range('E12').Value = '  +3.11%  '
Set-TextValue 'D13' '0.6506'
$ws.Range('E13').Value = '  +9.63%  '
Set-TextValue 'D14' '84.47'
$ws.Range('E14').Value = '  +10.82%  '
Set-TextValue 'D15' '4.733'
$ws.Range('E15').Value = '  +7.61%  '
Set-TextValue 'D16' '29.404.62'
$ws.Range('E16').Value = '  +10.61%  '
Set-TextValue 'D17' '0.9968'
$ws.Range('E17').Value = '  -0.49%  '
Set-TextValue 'D18' '0.000007343'
$ws.Range('E18').Value = '  +4.76%  '
Set-TextValue 'D19' '0.9965'
$ws.Range('E19').Value = '  -0.40%  '
Set-TextValue 'D20' '12.36'
$ws.Range('E20').Value = '  +10.14%  '
Set-TextValue 'D21' '2.061.57'
$ws.Range('E21').Value = '  +6.15%  '
Set-TextValue 'D22' '4.580'
$ws.Range('E22').Value = '  +5.19%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D23' '5.436'
$ws.Range('E23').Value = '  +7.86%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D24' '8.874'
$ws.Range('E24').Value = '  +6.93%  '
Set-TextValue 'D25' '143.84'
$ws.Range('E25').Value = '  +2.47%  '
Set-TextValue 'D26' '131.83'
$ws.Range('E26').Value = '  +25.16%  '
Set-TextValue 'D27' '16.49'
$ws.Range('E27').Value = '  +9.49%  '
Set-TextValue 'D28' '1.907'
$ws.Range('E28').Value = '  +7.40%  '
$ws.Range('E29').Value = '  +2.22%  '
Set-TextValue 'D30' '4.157'
$ws.Range('E30').Value = '  +5.80%  '
Set-TextValue 'D31' '0.08372'
$ws.Range('E31').Value = '  +6.67%  '
Set-TextValue 'D32' '3.798'
$ws.Range('E32').Value = '  +5.27%  '
Set-TextValue 'D33' '0.04943'
$ws.Range('E33').Value = '  +9.37%  '
Set-TextValue 'D34' '1.106'
$ws.Range('E34').Value = '  +11.97%  '
Set-TextValue 'D35' '0.6738'
$ws.Range('E35').Value = '  +10.78%  '
$ws.Range('E36').Value = '  +3.75%  '
Set-TextValue 'D37' '2.276'
$ws.Range('E37').Value = '  +16.71%  '
Set-TextValue 'D38' '2.716'
$ws.Range('E38').Value = '  +9.52%  '
Set-TextValue 'D39' '0.9549'
$ws.Range('E39').Value = '  +4.25%  '
Set-TextValue 'D40' '6.220'
$ws.Range('E40').Value = '  +9.63%  '
Set-TextValue 'D41' '0.01594'
$ws.Range('E41').Value = '  +7.98%  '
Set-TextValue 'D42' '0.9973'
$ws.Range('E42').Value = '  -0.27%  '
Set-TextValue 'D43' '102.33'
$ws.Range('E43').Value = '  +2.17%  '
Set-TextValue 'D44' '0.4087'
$ws.Range('E44').Value = '  +8.36%  '
Set-TextValue 'D45' '7.236'
$ws.Range('E45').Value = '  +8.19%  '
Set-TextValue 'D46' '0.1224'
$ws.Range('E46').Value = '  +7.40%  '
Set-TextValue 'D47' '0.05571'
$ws.Range('E47').Value = '  +4.25%  '
Set-TextValue 'D48' '31.89'
$ws.Range('E48').Value = '  +8.10%  '
Set-TextValue 'D49' '8.096'
$ws.Range('E49').Value = '  +4.34%  '
Set-TextValue 'D50' '1.314'
$ws.Range('E50').Value = '  +7.10%  '
Set-TextValue 'D51' '0.3627'
$ws.Range('E51').Value = '  +9.78%  '
